# Auto-generated edit script: update crypto price (D) and 1h volume change (E) columns
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "56.925.91"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -0.90%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.086.27"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +0.47%  "
$ws.Range("E4").Value = "  +0.06%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "521.74"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.17%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "136.82"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -2.42%  "
$ws.Range("E7").Value = "  +0.10%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "3.086.70"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +0.53%  "
$ws.Range("E9").Value = "  +2.91%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "7.36"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +2.92%  "
$ws.Range("E11").Value = "  -1.13%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.399"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +2.95%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "3.619.94"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +0.57%  "
$ws.Range("E14").Value = "  +1.50%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "25.49"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +0.71%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.0000161"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -1.06%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "57.085.46"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -0.73%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "3.091.61"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.75%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "5.89"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -2.71%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "12.48"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -0.95%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "7.88"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.33%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "347.81"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +2.85%  "
$ws.Range("E23").Value = "  +1.57%  "
$ws.Range("E24").Value = "  -0.05%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "68.00"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +1.78%  "
$ws.Range("E26").Value = "  -1.70%  "
$ws.Range("E27").Value = "  -0.33%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.996"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.37%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.0₃0881"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -2.74%  "
$ws.Range("E30").Value = "  -0.04%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "7.28"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +1.45%  "
$ws.Range("E32").Value = "  +0.69%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.90"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -6.85%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "20.77"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -0.15%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "4.96"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +8.44%  "
$ws.Range("E36").Value = "  -2.96%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "159.31"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +0.82%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "6.02"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -1.49%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "26.01"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +0.97%  "
$ws.Range("E40").Value = "  -0.34%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.0654"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -1.19%  "
$ws.Range("E42").Value = "  +0.49%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "4.03"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +1.77%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.693"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +1.89%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.387.39"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +5.26%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "36.59"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.35%  "
$ws.Range("E47").Value = "  +0.13%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "3.127.40"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +0.56%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0263"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +1.30%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.960"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -2.61%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "5.95"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -1.63%  "
